$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value changes: [cellRef, newValue, forceText]
$changes = @(
    @('D2', '62.723.43', $false),
    @('E2', '  +1.41%  ', $false),
    @('D3', '2.440.50', $false),
    @('E4', '  +0.11%  ', $false),
    @('D5', '566.79', $true),
    @('E5', '  +1.22%  ', $false),
    @('D6', '145.64', $true),
    @('E6', '  +2.43%  ', $false),
    @('E7', '  -0.10%  ', $false),
    @('E8', '  +0.22%  ', $false),
    @('E9', '  +2.46%  ', $false),
    @('D10', '0.154', $true),
    @('E10', '  +0.40%  ', $false),
    @('D11', '5.28', $true),
    @('E11', '  +0.76%  ', $false),
    @('D12', '0.355', $true),
    @('E12', '  +1.91%  ', $false),
    @('D13', '26.87', $true),
    @('E13', '  +5.50%  ', $false),
    @('E14', '  +5.16%  ', $false),
    @('D15', '2.882.04', $false),
    @('E15', '  +1.76%  ', $false),
    @('D16', '62.527.81', $false),
    @('E16', '  +1.16%  ', $false),
    @('D17', '2.442.97', $false),
    @('E17', '  +1.60%  ', $false),
    @('D18', '11.23', $true),
    @('E18', '  +0.45%  ', $false),
    @('E19', '  +2.27%  ', $false),
    @('D20', '323.82', $true),
    @('E20', '  +1.09%  ', $false),
    @('D21', '4.16', $true),
    @('E21', '  +0.90%  ', $false),
    @('D22', '0.999', $true),
    @('E22', '  -0.02%  ', $false),
    @('D23', '1.84', $true),
    @('E23', '  +6.43%  ', $false),
    @('D24', '67.32', $true),
    @('E24', '  +2.75%  ', $false),
    @('D25', '8.61', $true),
    @('E25', '  -1.76%  ', $false),
    @('D26', '583.61', $true),
    @('E26', '  +4.11%  ', $false),
    @('E27', '  +8.76%  ', $false),
    @('D28', '2.561.16', $false),
    @('B29', 'Binance-PegBSC-USD', $false),
    @('C29', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', $false),
    @('D29', '0.999', $true),
    @('E29', '  +0.15%  ', $false),
    @('B30', 'InternetComputer(DFINITY)', $false),
    @('C30', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', $false),
    @('D30', '8.41', $true),
    @('E30', '  +3.08%  ', $false),
    @('E31', '  +4.10%  ', $false),
    @('E32', '  +0.33%  ', $false),
    @('E33', '  +1.17%  ', $false),
    @('E34', '  +1.93%  ', $false),
    @('D35', '4.84', $true),
    @('E35', '  +2.50%  ', $false),
    @('E36', '  -0.13%  ', $false),
    @('E37', '  +1.58%  ', $false),
    @('D38', '18.80', $true),
    @('E38', '  +1.56%  ', $false),
    @('D39', '5.39', $true),
    @('E39', '  -0.15%  ', $false),
    @('D40', '148.16', $true),
    @('E40', '  -2.81%  ', $false),
    @('E41', '  +1.77%  ', $false),
    @('E42', '  +0.19%  ', $false),
    @('D43', '2.43', $true),
    @('E43', '  +8.79%  ', $false),
    @('D44', '148.58', $true),
    @('E44', '  +0.97%  ', $false),
    @('D45', '3.67', $true),
    @('E45', '  +2.13%  ', $false),
    @('D46', '0.0535', $true),
    @('E46', '  +1.36%  ', $false),
    @('D47', '20.49', $true),
    @('E47', '  +3.81%  ', $false),
    @('E48', '  +2.85%  ', $false),
    @('E49', '  +3.27%  ', $false),
    @('D50', '0.0921', $true),
    @('E50', '  +0.60%  ', $false),
    @('E51', '  +4.28%  ', $false)
)

foreach ($chg in $changes) {
    $cellRef = $chg[0]
    $newValue = $chg[1]
    $forceText = $chg[2]
    $range = $ws.Range($cellRef)
    if ($forceText) {
        # Original value is stored as literal text (e.g. keeps trailing zero,
        # like "18.80"), so force the cell to Text format before assigning,
        # otherwise Excel auto-converts it to a number and the formatting is lost.
        $range.NumberFormat = "@"
    }
    $range.Value = $newValue
}

Write-Host "Applied $($changes.Count) cell changes"
